# Add "2022-Q4" sheet data (feat: add 2022-Q4 data)
#
# 1. Insert a brand-new worksheet named "2022-Q4" right after "总计" and
#    before "2022-Q3".
# 2. Populate it with the 2022-Q4 fund-holdings table.
# 3. Update the "总计" (summary) sheet: push the existing quarter rows down
#    by one and insert the new 2022-Q4 totals at the top; append the extra
#    trailing index row that the growing table now needs.

$wb = $excel.ActiveWorkbook

# Helper: force a cell to be stored as literal text (so things like
# leading-zero fund codes, or numeric-looking strings such as "8.18", keep
# their exact textual representation instead of being auto-coerced to a
# number), then drop the temporary "@" number-format override again so the
# cell is left with the plain/default style.
function Set-TextCell($sheet, $addr, $text) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q4" worksheet right after "总计"
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item(1)
$newWs = $wb.Worksheets.Add($null, $totalWs)
$newWs.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# Step 2: populate the new sheet's header + data rows
# ---------------------------------------------------------------------------
$headerStyleSource = $totalWs.Range("B1")

$headers = @{
    "B" = "基金代码"
    "C" = "基金名称"
    "D" = "基金规模"
    "E" = "股票总仓位"
    "F" = "仓位占比"
    "G" = "持有市值(亿元)"
    "H" = "仓位排名"
}
foreach ($col in @("B", "C", "D", "E", "F", "G", "H")) {
    $headerAddr = "${col}1"
    $target = $newWs.Range($headerAddr)
    $target.Value = $headers[$col]
    $headerStyleSource.Copy()
    $target.PasteSpecial(-4122)  # xlPasteFormats
}

$rows = @(
    @{ A = 0;  B = "011335"; C = "银河医药健康混合A";             D = "8.18";  E = "92.65"; F = "4.60"; G = "0.3763"; H = 7  },
    @{ A = 1;  B = "012260"; C = "广发睿明优质企业混合A";          D = "10.10"; E = "65.70"; F = "2.69"; G = "0.2717"; H = 10 },
    @{ A = 2;  B = "000780"; C = "鹏华医疗保健股票";               D = "7.13";  E = "82.95"; F = "3.11"; G = "0.2217"; H = 8  },
    @{ A = 3;  B = "519673"; C = "银河康乐股票A";                 D = "2.15";  E = "93.79"; F = "5.23"; G = "0.1124"; H = 5  },
    @{ A = 4;  B = "002515"; C = "招商丰益灵活配置混合C";          D = "1.62";  E = "41.84"; F = "2.65"; G = "0.0429"; H = 8  },
    @{ A = 5;  B = "012261"; C = "广发睿明优质企业混合C";          D = "0.68";  E = "65.70"; F = "2.69"; G = "0.0183"; H = 10 },
    @{ A = 6;  B = "002514"; C = "招商丰益灵活配置混合A";          D = "0.60";  E = "41.84"; F = "2.65"; G = "0.0159"; H = 8  },
    @{ A = 7;  B = "010503"; C = "招商稳兴混合A";                 D = "0.93";  E = "28.39"; F = "1.63"; G = "0.0152"; H = 8  },
    @{ A = 8;  B = "009170"; C = "湘财长兴灵活配置混合C";          D = "0.37";  E = "87.98"; F = "3.13"; G = "0.0116"; H = 10 },
    @{ A = 9;  B = "009169"; C = "湘财长兴灵活配置混合A";          D = "0.18";  E = "87.98"; F = "3.13"; G = "0.0056"; H = 10 },
    @{ A = 10; B = "016018"; C = "银河康乐股票C";                 D = "0.10";  E = "93.79"; F = "5.23"; G = "0.0052"; H = 5  },
    @{ A = 11; B = "014692"; C = "中加量化研选混合型证券投资基金C"; D = "0.18";  E = "74.39"; F = "1.35"; G = "0.0024"; H = 5  },
    @{ A = 12; B = "015666"; C = "银河医药健康混合C";             D = "0.02";  E = "92.65"; F = "4.60"; G = "0.0009"; H = 7  },
    @{ A = 13; B = "014691"; C = "中加量化研选混合型证券投资基金A"; D = "0.04";  E = "74.39"; F = "1.35"; G = "0.0005"; H = 5  },
    @{ A = 14; B = "010504"; C = "招商稳兴混合C";                 D = "0.00";  E = "28.39"; F = "1.63"; G = 0;        H = 8  }
)

$rowNum = 2
foreach ($rowData in $rows) {
    $addrA = "A$rowNum"
    $addrB = "B$rowNum"
    $addrC = "C$rowNum"
    $addrD = "D$rowNum"
    $addrE = "E$rowNum"
    $addrF = "F$rowNum"
    $addrG = "G$rowNum"
    $addrH = "H$rowNum"

    $aCell = $newWs.Range($addrA)
    $aCell.Value = $rowData.A
    $headerStyleSource.Copy()
    $aCell.PasteSpecial(-4122)  # xlPasteFormats (matches the A-column index style)

    Set-TextCell $newWs $addrB $rowData.B
    Set-TextCell $newWs $addrC $rowData.C
    Set-TextCell $newWs $addrD $rowData.D
    Set-TextCell $newWs $addrE $rowData.E
    Set-TextCell $newWs $addrF $rowData.F
    if ($rowData.G -is [string]) {
        Set-TextCell $newWs $addrG $rowData.G
    } else {
        $newWs.Range($addrG).Value = $rowData.G
    }
    $newWs.Range($addrH).Value = $rowData.H

    $rowNum++
}

# ---------------------------------------------------------------------------
# Step 3: update the "总计" summary sheet
# ---------------------------------------------------------------------------
# Push the existing quarter rows (B2:D8) down to (B3:D9) - column A holds a
# plain positional index (0,1,2,...) and does not need to shift.
$totalWs.Range("B2:D8").Copy()
$totalWs.Range("B3").PasteSpecial(-4104)  # xlPasteAll

# The table grew by one row; give the new trailing row (row 9) the same
# style the other index cells in column A use.
$totalWs.Range("A8").Copy()
$totalWs.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$totalWs.Range("A9").Value = 7

# Write the new 2022-Q4 totals into the now-empty top data row.
$totalWs.Range("B2").Value = "2022-Q4"
$totalWs.Range("C2").Value = 15
$totalWs.Range("D2").Value = 1.1

$excel.CutCopyMode = 0
Write-Output "2022-Q4 sheet inserted and totals updated"
